$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A64").Value = "2025-04-29 08:55:35"
$ws.Range("B64").Value = 199
